$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "'5465859685115260427"
$ws.Range("I2").ClearFormats()
$ws.Range("I3").Value = "'5527442777154739395"
$ws.Range("I3").ClearFormats()
$ws.Range("I4").Value = "'11048837401662965790"
$ws.Range("I4").ClearFormats()
$ws.Range("I5").Value = "'17311461800651806091"
$ws.Range("I5").ClearFormats()
$ws.Range("I6").Value = "'10171788987322014488"
$ws.Range("I6").ClearFormats()
$ws.Range("I7").Value = "'1640453089486844230"
$ws.Range("I7").ClearFormats()
$ws.Range("I8").Value = "'15186096852338279725"
$ws.Range("I8").ClearFormats()
$ws.Range("I9").Value = "'18437259150481094143"
$ws.Range("I9").ClearFormats()
$ws.Range("I10").Value = "'1975812593444422882"
$ws.Range("I10").ClearFormats()
$ws.Range("I11").Value = "'17220968564977943993"
$ws.Range("I11").ClearFormats()
$ws.Range("I12").Value = "'9851901183800940666"
$ws.Range("I12").ClearFormats()
$ws.Range("I13").Value = "'4651596805960560616"
$ws.Range("I13").ClearFormats()
$ws.Range("I14").Value = "'5363012115113112816"
$ws.Range("I14").ClearFormats()
$ws.Range("I15").Value = "'12127878054617990148"
$ws.Range("I15").ClearFormats()
$ws.Range("I16").Value = "'12615781559068715528"
$ws.Range("I16").ClearFormats()
$ws.Range("I17").Value = "'16379321365746720216"
$ws.Range("I17").ClearFormats()
$ws.Range("I18").Value = "'1859052510171157074"
$ws.Range("I18").ClearFormats()
$ws.Range("I19").Value = "'131388509170113531"
$ws.Range("I19").ClearFormats()
$ws.Range("I20").Value = "'12276360953888945774"
$ws.Range("I20").ClearFormats()
$ws.Range("I21").Value = "'17948212404227027011"
$ws.Range("I21").ClearFormats()
$ws.Range("I22").Value = "'1777086285660115289"
$ws.Range("I22").ClearFormats()
$ws.Range("I23").Value = "'17550978780941298109"
$ws.Range("I23").ClearFormats()
$ws.Range("I24").Value = "'2131671840104369538"
$ws.Range("I24").ClearFormats()
$ws.Range("I25").Value = "'646713366152604332"
$ws.Range("I25").ClearFormats()
$ws.Range("I26").Value = "'10373905228459429881"
$ws.Range("I26").ClearFormats()
$ws.Range("I27").Value = "'4258573009579735995"
$ws.Range("I27").ClearFormats()
$ws.Range("I28").Value = "'265240268283751376"
$ws.Range("I28").ClearFormats()
$ws.Range("I29").Value = "'1716153287694692260"
$ws.Range("I29").ClearFormats()
$ws.Range("I30").Value = "'7011195061182511990"
$ws.Range("I30").ClearFormats()
$ws.Range("I31").Value = "'9986674937439859568"
$ws.Range("I31").ClearFormats()
$ws.Range("I32").Value = "'8924521215966159335"
$ws.Range("I32").ClearFormats()
$ws.Range("I33").Value = "'6389173515177744956"
$ws.Range("I33").ClearFormats()
$ws.Range("I34").Value = "'14148993971580128743"
$ws.Range("I34").ClearFormats()
$ws.Range("I35").Value = "'583765279824930520"
$ws.Range("I35").ClearFormats()
$ws.Range("I36").Value = "'10837231579467459848"
$ws.Range("I36").ClearFormats()
$ws.Range("I37").Value = "'10363795271348284161"
$ws.Range("I37").ClearFormats()
$ws.Range("I38").Value = "'15523820850506490992"
$ws.Range("I38").ClearFormats()
$ws.Range("I39").Value = "'5227751510528346188"
$ws.Range("I39").ClearFormats()
$ws.Range("I40").Value = "'7485263210229603469"
$ws.Range("I40").ClearFormats()
$ws.Range("I41").Value = "'15481276480988585061"
$ws.Range("I41").ClearFormats()
$ws.Range("I42").Value = "'11770231653918674322"
$ws.Range("I42").ClearFormats()
$ws.Range("I43").Value = "'13819626741178253463"
$ws.Range("I43").ClearFormats()
$ws.Range("I44").Value = "'13396730022424409992"
$ws.Range("I44").ClearFormats()
$ws.Range("I45").Value = "'2523881632023293484"
$ws.Range("I45").ClearFormats()
$ws.Range("I46").Value = "'10587417906698095327"
$ws.Range("I46").ClearFormats()
$ws.Range("I47").Value = "'9033843866905793904"
$ws.Range("I47").ClearFormats()
$ws.Range("I48").Value = "'8368446460655908217"
$ws.Range("I48").ClearFormats()
$ws.Range("I49").Value = "'2448008120671130484"
$ws.Range("I49").ClearFormats()
$ws.Range("I50").Value = "'676658871738576226"
$ws.Range("I50").ClearFormats()
$ws.Range("I51").Value = "'11815689254016262160"
$ws.Range("I51").ClearFormats()
$ws.Range("I52").Value = "'1865281822864486586"
$ws.Range("I52").ClearFormats()
$ws.Range("I53").Value = "'11999168590972673814"
$ws.Range("I53").ClearFormats()
$ws.Range("I54").Value = "'15663199238863036415"
$ws.Range("I54").ClearFormats()
$ws.Range("I55").Value = "'3482048972607296362"
$ws.Range("I55").ClearFormats()
$ws.Range("I56").Value = "'9297669140363783287"
$ws.Range("I56").ClearFormats()
$ws.Range("I57").Value = "'825577570103515494"
$ws.Range("I57").ClearFormats()
$ws.Range("I58").Value = "'17077985264131744185"
$ws.Range("I58").ClearFormats()
$ws.Range("I59").Value = "'2590433911241102872"
$ws.Range("I59").ClearFormats()
$ws.Range("I60").Value = "'583355402538321026"
$ws.Range("I60").ClearFormats()
$ws.Range("I61").Value = "'389810736453788751"
$ws.Range("I61").ClearFormats()
$ws.Range("I62").Value = "'7050301517988511742"
$ws.Range("I62").ClearFormats()
$ws.Range("I63").Value = "'12641924004534609258"
$ws.Range("I63").ClearFormats()
$ws.Range("I64").Value = "'9056830240370997859"
$ws.Range("I64").ClearFormats()
$ws.Range("I65").Value = "'13942523387772293193"
$ws.Range("I65").ClearFormats()
$ws.Range("I66").Value = "'15832826006729248969"
$ws.Range("I66").ClearFormats()
$ws.Range("I67").Value = "'13355582891255512253"
$ws.Range("I67").ClearFormats()
$ws.Range("I68").Value = "'17377645535285548270"
$ws.Range("I68").ClearFormats()
$ws.Range("I69").Value = "'6897379166409132309"
$ws.Range("I69").ClearFormats()
$ws.Range("I70").Value = "'16292580299745028902"
$ws.Range("I70").ClearFormats()
$ws.Range("I71").Value = "'4397641155899062092"
$ws.Range("I71").ClearFormats()
